# Append the new job-book row (row 62) to the JOBS sheet, mirroring the
# existing rows' shape: jobNumber is numeric, everything else is text
# (including the boolean "_isDeleted" column), matching the diff that adds
# job 71314 / DXS36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

# --- jobNumber (numeric) ---------------------------------------------------
$ws.Cells.Item($row, 1).Value = 71314

# --- modelNumber / serialNumber (plain text) -------------------------------
$ws.Cells.Item($row, 2).Value = "DXS36"
$ws.Cells.Item($row, 3).Value = "?"

# --- voltage: looks numeric ("460") but must stay text, like the rest of
#     the sheet (column is full of numberStoredAsText values) -------------
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "460"
$ws.Cells.Item($row, 4).Style = "Normal"

# --- unloaders / statorStatus ----------------------------------------------
$ws.Cells.Item($row, 5).Value = "?"
$ws.Cells.Item($row, 6).Value = "?"

# --- incomingNumber: also numeric-looking text ("82863") -------------------
$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "82863"
$ws.Cells.Item($row, 7).Style = "Normal"

# --- scrap -------------------------------------------------------------
$ws.Cells.Item($row, 8).Value = "NO"

# --- notes: empty text ---------------------------------------------------
$ws.Cells.Item($row, 9).Value = ""

# --- enteredBy -------------------------------------------------------------
$ws.Cells.Item($row, 10).Value = "wes"

# --- enteredOn: date-looking text ("9/12/2022") - keep as literal text -----
$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = "9/12/2022"
$ws.Cells.Item($row, 11).Style = "Normal"

# --- _isDeleted: boolean FALSE ----------------------------------------------
$ws.Cells.Item($row, 12).Value = $false

# --- deletedBy / deletedOn --------------------------------------------------
$ws.Cells.Item($row, 13).Value = "N/A"
$ws.Cells.Item($row, 14).Value = "N/A"

# --- warranty ----------------------------------------------------------
$ws.Cells.Item($row, 15).Value = "NO"
